# Edit script: "Quantum Entanglement" essay -> "Mathematics" essay
# Applies the textual substitutions described by the target diff.
#
# Note: we locate text with Range.Find.Execute (search-only, no Replace
# argument) and then assign Range.Text directly. Going through
# Find.Execute's own Replace mechanism triggers AutoCorrect "smart quotes",
# which the target text does not use (it keeps plain straight apostrophes).

$d = $word.ActiveDocument

function Set-FoundText($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# Appends `extra` right after the (single) existing occurrence of `anchor`,
# as new, separate run(s) of text - used for the two spots where the target
# grows a sentence into two sentences.
function Append-AfterText($anchor, $extra) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $ok) {
        Write-Output "NOT FOUND (append): $anchor"
        return
    }
    $rng.Collapse(0) | Out-Null
    $rng.InsertAfter($extra)
}

# --- Title ---
Set-FoundText "Quantum Entanglement: Unveiling the Mysteries of Interconnectedness" "Mathematics: The Universal Language of Science and Logic"

# --- Author name ---
Set-FoundText " Samuel Davies" " Jasper Lancaster"

# --- Author email/handle (collapses the 5 split runs into one run of new text) ---
Set-FoundText "samuel.davies@quantumstudies.org" "yourvalidname"

# --- Body paragraph 1 (first block, before the "In 1935..." block) ---
Set-FoundText "Amidst the perplexing wonders of the quantum realm exists a profound phenomenon challenging our understanding of reality - quantum entanglement" "Mathematics, a subject that has fascinated and challenged minds for centuries, unveils the mysteries of the universe through the intricate tapestry of numbers, equations, and geometric patterns"

Set-FoundText " This enigmatic connection between particles, regardless of their distance, has captured the imagination of scientists, philosophers, and artists alike" " Mathematics provides a universal language that transcends cultures, enabling us to understand the cosmos, unravel the enigmas of nature, and harness the power of logic to solve complex problems"

Set-FoundText " In this essay, we embark on a journey to unravel the mysteries of quantum entanglement, exploring its implications for our comprehension of the universe and delving into the potential applications that may revolutionize various fields" " The field of mathematics is a symphony of abstract concepts, where symbols dance in harmony, revealing the underlying order and beauty of our world"

# --- Body paragraph 1, second block ("In 1935..." / "Their proposal..." / "This relationship...") ---
Set-FoundText "In 1935, Albert Einstein, Boris Podolsky, and Nathan Rosen introduced the concept of quantum entanglement through their famous thought experiment known as the EPR paradox" "Through the exploration of mathematical concepts, we unlock the secrets of nature's blueprint"

Set-FoundText " Their proposal demonstrated that two particles, once entangled, remain interconnected regardless of the distance separating them" " The Fibonacci sequence, found in the spirals of seashells or the patterns of plant growth, exemplifies the intricate relationship between numbers and biological structures"

# This sentence gains two follow-on sentences ("." + "These patterns underscore...") after it in the target.
Set-FoundText " This relationship transcends the constraints of space and time, allowing one particle to instantaneously influence the other, even across vast cosmological distances" " The elegance of geometric shapes, such as fractals, reflects the self-similarity found in everything from snowflakes to coastlines"

Append-AfterText " The elegance of geometric shapes, such as fractals, reflects the self-similarity found in everything from snowflakes to coastlines" ". These patterns underscore the profound interconnectedness of all things and provide a glimpse into the underlying mathematical principles that govern our universe"

# --- Body paragraph 1, third block ("Moreover..." / "When entangled..." / "This non-locality...") ---
Set-FoundText "Moreover, quantum entanglement defies classical intuition" "Mathematics isn't merely a collection of abstract theories; it's a powerful tool with practical applications in every field imaginable"

Set-FoundText " When entangled particles are measured, their properties, such as spin or polarization, are correlated in a way that cannot be explained by classical physics" " It empowers engineers to design structures that withstand earthquakes, enables us to predict weather patterns, and makes it possible to develop new medical treatments and technologies"

# This sentence gains two follow-on sentences ("." + "It's a subject that touches...") after it in the target.
Set-FoundText " This non-locality, as it is known, challenges our conventional notions of causality and raises fundamental questions about the nature of reality itself" " From the economy to finance, from computer science to data analysis, and even in music and art, the profound influence of mathematics is undeniable"

$rng = $d.Content
$rng.Find.Execute(" From the economy to finance, from computer science to data analysis, and even in music and art, the profound influence of mathematics is undeniable", $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(". It's a subject that touches every aspect of our lives, shaping our understanding of the world and guiding us towards a future filled with infinite possibilities")

# --- Summary paragraph ---
Set-FoundText "Quantum entanglement, an awe-inspiring phenomenon, offers a glimpse into the uncharted territory of the quantum world" "Mathematics serves as an essential tool for understanding the intricacies of our universe, providing a lens through which we unlock the mysteries of science and logic"

Set-FoundText " Its non-local nature challenges our fundamental understanding of reality, while its potential applications hold promise for transformative technologies" " Its abstract concepts find practical applications in diverse fields, empowering engineers, scientists, musicians, artists, and countless other professionals to innovate and drive progress"

Set-FoundText " From quantum computing to secure communication, entanglement-based technologies may revolutionize numerous fields" " Mathematics unveils the interconnectedness of all things, from the Fibonacci sequence found in nature to the intricate patterns in art and music"

Set-FoundText " Though much remains unknown, continued exploration of quantum entanglement promises to deepen our comprehension of the universe and expand the boundaries of human knowledge" " It's a subject that permeates our existence, shaping our understanding of the cosmos and enabling us to chart a course toward a future of endless possibilities"

# --- Add a trailing empty paragraph at the end of the document ---
$d.Content.InsertParagraphAfter() | Out-Null
